$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Text / date-string fields (force text storage so e.g. "001" keeps its
# leading zero instead of being coerced to the number 1)
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "001"
$ws.Range("M2").NumberFormat = "@"
$ws.Range("M2").Value = "2020-12-16 00:00:00"
$ws.Range("N2").NumberFormat = "@"
$ws.Range("N2").Value = "2019-12-31 00:00:00"

# Numeric fields
$ws.Range("O2").Value = 170787524.38
$ws.Range("P2").Value = 1654020988.88
$ws.Range("Q2").Value = 1484836060.35
$ws.Range("R2").Value = -17.1411591581
$ws.Range("S2").Value = 1352291071.35
$ws.Range("T2").Value = 1352291071.35
$ws.Range("U2").Value = -17.4689622231
$ws.Range("V2").Value = 37948854.87
$ws.Range("W2").Value = 20883482.3
$ws.Range("X2").Value = 19468483.32
$ws.Range("Y2").Value = 194180217.99
$ws.Range("Z2").Value = 193572092.06
$ws.Range("AA2").Value = 22732789.55
$ws.Range("AG2").Value = 5762868.3
$ws.Range("AP2").Value = -14.0562614823
$ws.Range("AQ2").Value = 41.930961499816
$ws.Range("AR2").Value = 36.245369994984
$ws.Range("AS2").Value = 165255024.38
$ws.Range("AT2").Value = 34.720908335118
